$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 631 (shifts existing rows 631.. down by one)
$ws.Rows.Item(631).Insert()

# Populate the newly inserted row 631 with the new data record
$ws.Range("A631").Value = 4
$ws.Range("B631").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C631").Value = "Los Lagos"
$ws.Range("D631").Value = 45041
$ws.Range("E631").Value = 10
$ws.Range("F631").Value = 100114001
$ws.Range("G631").Value = "Papa"
$ws.Range("H631").Value = "Patagonia"
$ws.Range("I631").Value = "1a (cosecha)"
$ws.Range("J631").Value = 600
$ws.Range("K631").Value = 12000
$ws.Range("L631").Value = 12000
$ws.Range("M631").Value = 12000
$ws.Range("N631").Value = "$/saco 25 kilos"
$ws.Range("O631").Value = "Provincia de Llanquihue"
$ws.Range("P631").Value = 480
$ws.Range("Q631").Value = 25
$ws.Range("R631").Value = "Hortaliza"
